$wb = $excel.ActiveWorkbook

# --- Sheet "pro" (sheet1): update B2:B26, column width, selection ---
$wsPro = $wb.Worksheets.Item("pro")
$proValues = @(3028696.6213632673,3273601.8239796101,3536322.1498745177,3361751.4325653021,3553122.3174919696,3591999.7167676925,3297147.2015583459,3033763.8545851018,3348183.3025720734,3425947.6561750858,3412867.593579961,3489372.1905571036,3293582.3147831326,3431439.3985320223,3711735.7454813719,3617139.0391392587,4075650.3465787941,4565920.1597745055,4764012.0771354465,4703192,4847703,4882074.1284464085,5165032.3033993607,5452393.3539396049,5166418.8881907323)
for ($i = 0; $i -lt $proValues.Length; $i++) {
    $wsPro.Cells.Item($i + 2, 2).Value = $proValues[$i]
}
$wsPro.Columns.Item(2).ColumnWidth = 14.6328125

# --- Sheet "ind" (sheet2): update B2:B101, selection ---
$wsInd = $wb.Worksheets.Item("ind")
$indValues = @(699670.79552178376,647354.05799986003,612780.99896056426,655825.32636813412,737881.65011925215,757940.66238015192,681940.23304213362,759595.43208880175,827400.3465482546,860063.38071826799,759976.78665241238,843613.0228515428,840822.31895228114,878115.38618673733,749827.27320059587,785720.41148036229,824762.82411476399,853027.20948283549,867600.21700046002,758206.90324583906,782330.29076591053,896267.91784728167,850794.12563381926,810351.97528692498,816339.58560376859,805501.69538828766,707540.65135222918,518965.61597861868,543735.53916052356,531667.47230982874,536430.85443877522,572495.63944459846,503505.00693830429,537922.32139363373,662032.43980060739,534545.85808071657,609787.34854261565,592986.8435967617,562784.81643530983,533618.31013936643,537926.00370099372,548658.96180714283,552565.42662041157,591934.88118156511,640361.1975835989,588411.72586998018,658699.5915981459,568646.75989242143,542466.41741161817,514124.09238478902,598694.68923233123,663020.37721375294,867899.62126419682,851188.12655503803,876765.13565364713,884359.57272253674,966930.81022490398,1019201.374594875,996706.18917375361,781654.44527532882,591408.85960995185,580786.89655719767,783379.07096350146,830283.31373213045,842785.14394586533,872370.83619518951,863963.52322285692,906176.04966184602,1020157.4211890965,1000885.0618747902,811581.93860865908,932937.66357818525,1024258.2843377531,1017745.8880120541,890180.51961162558,996745.86827503168,916791.07998724503,983518.98814131983,973677.67992797715,1004783.8587884611,1039781.4083080392,1049899.9601661121,971131.61115762126,1059291.8263736228,1174021.2760110451,1166578.3102370845,1090040.8385906636,1161961.5336492332,1239239.7474453468,1184285.1581796382,1092338.8112512447,1285387.0542468375,1311363.6644483048,1362878.4809335119,1228698.8059495161,1265461.6720702609,1326488.8555946348,1289393.1299960888,1263903.4225931147,1203261.8826738598)
for ($i = 0; $i -lt $indValues.Length; $i++) {
    $wsInd.Cells.Item($i + 2, 2).Value = $indValues[$i]
}

# --- Sheet "conso" (sheet4): update B2:B26, selection ---
$wsConso = $wb.Worksheets.Item("conso")
$consoValues = @(1942740.8761631781,2095452.5588326475,2288116.6680541569,2159460.1663240478,2285711.9595365045,2325892.1426933324,2184620.5898217503,2007436.0413187987,2209565.3358346364,2264399.5858566375,2262810.7966729659,2312842.6791132903,2179077.0539652337,2278716.0130595719,2483983.5664652353,2435322.2488249354,2737230.4804706117,3061555.0096661635,3185425.4556566388,3144774,3388659,3200386.2455012747,3489000.1367686647,3683113.6845149263,3489936.780383375)
for ($i = 0; $i -lt $consoValues.Length; $i++) {
    $wsConso.Cells.Item($i + 2, 2).Value = $consoValues[$i]
}

# --- Sheet "VA" (sheet3) is formula-driven (=pro!Bn - conso!Bn); recalculates automatically ---
$wsVA = $wb.Worksheets.Item("VA")

# --- Update selections on every sheet to B2:B26 (matches the authored diff) ---
$wsInd.Activate()
$wsInd.Range("B2:B26").Select()

$wsVA.Activate()
$wsVA.Range("B2:B26").Select()

$wsConso.Activate()
$wsConso.Range("B2:B26").Select()

# "pro" stays the active/tab-selected sheet, selection B2:B26
$wsPro.Activate()
$wsPro.Range("B2:B26").Select()

$excel.Calculate()
